# Weekly update: insert the newest week's Mango price row at the top of the
# data block (row 110) and push the existing history rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 110:119 down to 111:120, inserting a fresh blank row 110.
$ws.Rows("110").Insert()

# Populate the new row 110 with this week's data.
$ws.Range("A110").Value = 11
$ws.Range("B110").Value = "Vega Monumental Concepción"
$ws.Range("C110").Value = "Bíobío"
$ws.Range("D110").Value = 44769
$ws.Range("E110").Value = 8
$ws.Range("F110").Value = "Fruta"
$ws.Range("G110").Value = 100108
$ws.Range("H110").Value = "Tropicales y subtropicales"
$ws.Range("I110").Value = 100108002
$ws.Range("J110").Value = "Mango"
$ws.Range("K110").Value = "Sin especificar"
$ws.Range("L110").Value = "Primera"
$ws.Range("M110").Value = 220
$ws.Range("N110").Value = 9000
$ws.Range("O110").Value = 9500
$ws.Range("P110").Value = 9227
$ws.Range("Q110").Value = "$/bandeja 4 kilos"
$ws.Range("R110").Value = "Perú"
$ws.Range("S110").Value = 2307
$ws.Range("T110").Value = 4
